# u matrix: electricity use from gas boiler put to zero
#
# Row 4 ("Electricity" need) had non-zero electricity allocations for
# activities that are actually gas-fired (gas boiler for heating, gas
# boiler for hot sanitary water, gas stove for cooking). Those spurious
# electricity contributions are zeroed out:
#   H4 (Exploiting Gas boiler for Heating)              0.1  -> 0
#   K4 (Exploiting Gas boiler for Hot Sanitary Water)   0.1  -> 0
#   P4 (Exploiting Gas Stove for Cooking)               0.01 -> 0

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("P4").Value = 0

# Reflect the author's final cursor position/selection on the sheet.
$ws.Range("N4").Select()
